$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Retrain results: column J -> 1, column K -> 0.5 for all data rows (1-51),
# replacing the old string-labeled header values and the previous 0.5/1 split.
$ws.Range("J1:J51").Value = 1
$ws.Range("K1:K51").Value = 0.5

# Update the view: scroll near the bottom of the data and select the full K column's data range.
$ws.Range("K1:K51").Select()
